$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Sample ID" query text (shared string originally used by cell B3)
# is being trimmed: the Tumor and Analyte Type columns are removed from
# the SELECT list. Re-assigning the full corrected text to B3 causes the
# workbook's shared-string table to be rebuilt, which is also why the
# "File Name" query (previously after it in the table) now sorts ahead
# of it - that is an automatic side effect of the edit below, not a
# separate change.
$newSampleIdQuery = @'
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
   s.phs_accession = 'phs001819' AND gi.reference_genome_assembly = 'GRCh37'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
'@

$ws.Range("B3").Value = $newSampleIdQuery

# Reflect that the user had scrolled down and was focused on the edited
# cell (B3) when the workbook was saved.
$ws.Range("B3").Select()
